# Scheduled-runner style refresh of the Leve profit-calc columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 13 (Leve Item ID 2144)
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 10000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -9831
$ws.Range("N13").ClearContents()

# Row 15 (Leve Item ID 44146)
$ws.Range("H15").Value = 480.13
$ws.Range("I15").Value = 480.13
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1440.39
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1271.39

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 2002.4
$ws.Range("I100").Value = 2000

# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 59685.36
$ws.Range("I113").Value = 91827.25
$ws.Range("J113").Value = 2544.2222
$ws.Range("K113").Value = 91827.25
$ws.Range("L113").Value = 2544.2222
$ws.Range("M113").Value = -88573.25
$ws.Range("N113").Value = -9052.2222

# Row 139 (Leve Item ID 42306)
$ws.Range("H139").Value = 49920
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49920
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49920
$ws.Range("N139").Value = -60200


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 8773170
$ws.Range("I61").Value = 10639551
$ws.Range("J61").Value = 1183
$ws.Range("K61").Value = 10639551
$ws.Range("L61").Value = 1183
$ws.Range("M61").Value = -10639339
$ws.Range("N61").Value = -1607

# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 3252.5
$ws.Range("I63").Value = 3503.3333
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 3503.3333
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -2817.3333
$ws.Range("N63").Value = -3872

# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 3252.5
$ws.Range("I66").Value = 3503.3333
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 17516.6665
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -14084.6665
$ws.Range("N66").Value = -19364

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 5557432
$ws.Range("I132").Value = 6580696.5
$ws.Range("J132").Value = 2565.1428
$ws.Range("K132").Value = 19742089.5
$ws.Range("L132").Value = 7695.428400000001
$ws.Range("M132").Value = -19739559.5
$ws.Range("N132").Value = -12755.4284

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 8773170
$ws.Range("I136").Value = 10639551
$ws.Range("J136").Value = 1183
$ws.Range("K136").Value = 31918653
$ws.Range("L136").Value = 3549
$ws.Range("M136").Value = -31916103
$ws.Range("N136").Value = -8649


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3600
$ws.Range("I134").Value = 2881.6128
$ws.Range("J134").Value = 5084.6665
$ws.Range("K134").Value = 8644.838400000001
$ws.Range("L134").Value = 15253.9995
$ws.Range("M134").Value = -6109.838400000001
$ws.Range("N134").Value = -20323.9995


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 106 (Leve Item ID 18661)
$ws.Range("H106").Value = 47396
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 47396
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 47396
$ws.Range("N106").Value = -49920

# Row 135 (Leve Item ID 42008)
$ws.Range("H135").Value = 54950
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 54950
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 54950
$ws.Range("N135").Value = -65090


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 4040156.8
$ws.Range("I4").Value = 4166830
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 12500490
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -12500378
$ws.Range("N4").Value = -3000224

# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 1558.2858
$ws.Range("I5").Value = 301.35715
$ws.Range("J5").Value = 4072.1428
$ws.Range("K5").Value = 904.0714499999999
$ws.Range("L5").Value = 12216.4284
$ws.Range("M5").Value = -792.0714499999999
$ws.Range("N5").Value = -12440.4284

# Row 25 (Leve Item ID 4709)
$ws.Range("H25").Value = 700.5
$ws.Range("I25").Value = 266.66666
$ws.Range("J25").Value = 2002
$ws.Range("K25").Value = 799.9999799999999
$ws.Range("L25").Value = 6006
$ws.Range("M25").Value = -630.9999799999999
$ws.Range("N25").Value = -6344

# Row 30 (Leve Item ID 4709)
$ws.Range("H30").Value = 700.5
$ws.Range("I30").Value = 266.66666
$ws.Range("J30").Value = 2002
$ws.Range("K30").Value = 799.9999799999999
$ws.Range("L30").Value = 6006
$ws.Range("M30").Value = -697.9999799999999
$ws.Range("N30").Value = -6210

# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1284
$ws.Range("I68").Value = 690
$ws.Range("J68").Value = 1383
$ws.Range("K68").Value = 2070
$ws.Range("L68").Value = 4149
$ws.Range("M68").Value = -1259
$ws.Range("N68").Value = -5771

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1284
$ws.Range("I71").Value = 690
$ws.Range("J71").Value = 1383
$ws.Range("K71").Value = 6210
$ws.Range("L71").Value = 12447
$ws.Range("M71").Value = -2154
$ws.Range("N71").Value = -20559

# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 1869.7273
$ws.Range("I80").Value = 288.5
$ws.Range("J80").Value = 2221.111
$ws.Range("K80").Value = 865.5
$ws.Range("L80").Value = 6663.333
$ws.Range("M80").Value = 70.5
$ws.Range("N80").Value = -8535.332999999999

# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 1869.7273
$ws.Range("I83").Value = 288.5
$ws.Range("J83").Value = 2221.111
$ws.Range("K83").Value = 2596.5
$ws.Range("L83").Value = 19989.999
$ws.Range("M83").Value = 2083.5
$ws.Range("N83").Value = -29349.999

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1977.5
$ws.Range("I132").Value = 720
$ws.Range("J132").Value = 3235
$ws.Range("K132").Value = 6480
$ws.Range("L132").Value = 29115
$ws.Range("M132").Value = -3950
$ws.Range("N132").Value = -34175

# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 1558.2858
$ws.Range("I135").Value = 301.35715
$ws.Range("J135").Value = 4072.1428
$ws.Range("K135").Value = 2712.21435
$ws.Range("L135").Value = 36649.2852
$ws.Range("M135").Value = -177.2143499999997
$ws.Range("N135").Value = -41719.2852


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 18 (Leve Item ID 4309)
$ws.Range("H18").Value = 6250
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 12000
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = -207
$ws.Range("N18").Value = -12586

# Row 43 (Leve Item ID 4218)
$ws.Range("H43").Value = 6227.143
$ws.Range("I43").Value = 1196.6666

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 3416.4614
$ws.Range("I102").Value = 4224.3335
$ws.Range("J102").Value = 1598.75
$ws.Range("K102").Value = 4224.3335
$ws.Range("L102").Value = 1598.75
$ws.Range("M102").Value = -2602.3335
$ws.Range("N102").Value = -4842.75

# Row 133 (Leve Item ID 41854)
$ws.Range("H133").Value = 50779
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50779
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50779
$ws.Range("N133").Value = -60899


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 4656.523
$ws.Range("I7").Value = 4065.8147
$ws.Range("J7").Value = 5594.706
$ws.Range("K7").Value = 4065.8147
$ws.Range("L7").Value = 5594.706
$ws.Range("M7").Value = -3953.8147
$ws.Range("N7").Value = -5818.706

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 4656.523
$ws.Range("I126").Value = 4065.8147
$ws.Range("J126").Value = 5594.706
$ws.Range("K126").Value = 12197.4441
$ws.Range("L126").Value = 16784.118
$ws.Range("M126").Value = -9727.444100000001
$ws.Range("N126").Value = -21724.118


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 19 (Leve Item ID 2666)
$ws.Range("H19").Value = 50000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 50000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 50000
$ws.Range("N19").Value = -50348

